# Trade #105 closed at 2026-02-16 21:41:02 - leadlag DOWN +0.000%
#
# This script applies the workbook edit described by the commit:
#   - Refreshes aggregate stats on "Summary" and "Comparison" sheets
#   - Closes out the previously-OPEN trade rows on "leadlag" (row 57) and
#     "momentum" (row 19), and mirrors those now-CLOSED rows onto "All Trades"
#   - Appends a brand-new OPEN trade (#105) to the "leadlag" sheet

$wb = $excel.ActiveWorkbook

# Helper: write a value as literal text (mirrors the source file's use of
# inline/shared strings for things that look numeric/date-like, e.g. "68.9%",
# "2026-02-16", "3.08") instead of letting Excel auto-coerce them into
# numbers/dates. An empty string needs special handling because plain
# assignment of "" clears the cell instead of leaving a zero-length Text cell.
function Set-Text($cell, [string]$val) {
    if ($val -eq "") {
        $cell.Value = "'"
    } else {
        $cell.NumberFormat = "@"
        $cell.Value = $val
    }
    $cell.Style = "Normal"
}

# ---------------------------------------------------------------------
# Summary sheet
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Summary")

$ws.Range("C2").Value = 74
Set-Text $ws.Range("D2") "68.9%"
Set-Text $ws.Range("E2") "+21.8819%"
Set-Text $ws.Range("F2") "+0.2957%"

$ws.Range("C3").Value = 79
Set-Text $ws.Range("D3") "44.3%"
Set-Text $ws.Range("E3") "+12.8142%"
Set-Text $ws.Range("F3") "+0.1622%"

Set-Text $ws.Range("D4") "64.0%"
Set-Text $ws.Range("E4") "+9.0677%"
Set-Text $ws.Range("F4") "+0.3627%"

# ---------------------------------------------------------------------
# leadlag sheet
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("leadlag")

# Row 57 (Trade #73) transitions from OPEN -> CLOSED
$ws.Range("G57").Value = 68723.02183500001
Set-Text $ws.Range("H57") "CLOSED"
$ws.Range("I57").Value = 0.1393
$ws.Range("J57").Value = 1.39
Set-Text $ws.Range("M57") "time_exit_5min"
$ws.Range("N57").Value = 5

# New row 81 - Trade #105, newly opened
$ws.Cells.Item(81, 1).Value = 105
Set-Text $ws.Cells.Item(81, 2) "2026-02-16"
Set-Text $ws.Cells.Item(81, 3) "21:41:02"
Set-Text $ws.Cells.Item(81, 4) "leadlag"
Set-Text $ws.Cells.Item(81, 5) "DOWN"
$ws.Cells.Item(81, 6).Value = 68395.58
Set-Text $ws.Cells.Item(81, 7) ""
Set-Text $ws.Cells.Item(81, 8) "OPEN"
$ws.Cells.Item(81, 9).Value = 0
$ws.Cells.Item(81, 10).Value = 0
$ws.Cells.Item(81, 11).Value = 0.75
Set-Text $ws.Cells.Item(81, 12) "Binance leading with -0.126% move"
Set-Text $ws.Cells.Item(81, 13) ""
$ws.Cells.Item(81, 14).Value = 0

# ---------------------------------------------------------------------
# momentum sheet
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("momentum")

# Row 19 (Trade #74) transitions from OPEN -> CLOSED
$ws.Range("G19").Value = 68876.249928
Set-Text $ws.Range("H19") "CLOSED"
$ws.Range("I19").Value = 0.2503
$ws.Range("J19").Value = 2.5
Set-Text $ws.Range("M19") "time_exit_5min"
$ws.Range("N19").Value = 5

# ---------------------------------------------------------------------
# All Trades sheet - mirror the two newly-closed trades
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("All Trades")

# New row 74 == leadlag row 57 (Trade #73), now CLOSED
$ws.Cells.Item(74, 1).Value = 73
Set-Text $ws.Cells.Item(74, 2) "2026-02-16"
Set-Text $ws.Cells.Item(74, 3) "21:35:50"
Set-Text $ws.Cells.Item(74, 4) "leadlag"
Set-Text $ws.Cells.Item(74, 5) "UP"
$ws.Cells.Item(74, 6).Value = 68627.41
$ws.Cells.Item(74, 7).Value = 68723.02183500001
Set-Text $ws.Cells.Item(74, 8) "CLOSED"
$ws.Cells.Item(74, 9).Value = 0.1393
$ws.Cells.Item(74, 10).Value = 1.39
$ws.Cells.Item(74, 11).Value = 0.75
Set-Text $ws.Cells.Item(74, 12) "Binance leading with 0.116% move"
Set-Text $ws.Cells.Item(74, 13) "time_exit_5min"
$ws.Cells.Item(74, 14).Value = 5

# New row 75 == momentum row 19 (Trade #74), now CLOSED
$ws.Cells.Item(75, 1).Value = 74
Set-Text $ws.Cells.Item(75, 2) "2026-02-16"
Set-Text $ws.Cells.Item(75, 3) "21:35:56"
Set-Text $ws.Cells.Item(75, 4) "momentum"
Set-Text $ws.Cells.Item(75, 5) "UP"
$ws.Cells.Item(75, 6).Value = 68704.28
$ws.Cells.Item(75, 7).Value = 68876.249928
Set-Text $ws.Cells.Item(75, 8) "CLOSED"
$ws.Cells.Item(75, 9).Value = 0.2503
$ws.Cells.Item(75, 10).Value = 2.5
$ws.Cells.Item(75, 11).Value = 0.9
Set-Text $ws.Cells.Item(75, 12) "Upward momentum: 0.183% over 10 samples"
Set-Text $ws.Cells.Item(75, 13) "time_exit_5min"
$ws.Cells.Item(75, 14).Value = 5

# ---------------------------------------------------------------------
# Comparison sheet
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Comparison")

$ws.Range("B2").Value = 79
Set-Text $ws.Range("C2") "44.3%"
Set-Text $ws.Range("D2") "3.08"
Set-Text $ws.Range("E2") "+0.5421%"
Set-Text $ws.Range("G2") "1.85"

Set-Text $ws.Range("C3") "64.0%"
Set-Text $ws.Range("D3") "9.06"
Set-Text $ws.Range("E3") "+0.6370%"
Set-Text $ws.Range("G3") "1.13"

Write-Output "edit applied"
